# Add a new "U2 / Average Wind Speed / m/s" row to the weather station
# metadata table, inserted just above the existing "UMax" row (i.e. it
# becomes the new row 18, pushing the old row 18 down to row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 18 (the current "UMax" row), shifting
# that row and everything below it down by one.
$ws.Rows("18:18").Insert()

# Copy the formatting (style) from the row that got pushed down to row 19
# onto the freshly inserted blank row 18, so the new row matches the look
# of the rest of the table.
$ws.Range("A19:C19").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row's content.
$ws.Range("A18").Value = "U2"
$ws.Range("B18").Value = "Average Wind Speed"
$ws.Range("C18").Value = "m/s"

# Reflect the author's final cell selection in the saved file.
$ws.Range("C10").Select()
